$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10 (columns A-T), following Dr Hou advice: Sending/Target clusters
# now cover all combinations of ECs/FAPs/sCs, with updated statistics.
$rowData = @(
    @("ECs", "Efnb1", "Ephb2", "ECs", 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 1, 0.3333333333333333, 0.005966, 0.017898, 0.00125947234315407, 0.00125947234315407, 0.02324652894333333, 0.20921876049, 0.0005059558618197952, 0.0005059558618197952),
    @("ECs", "Efnb1", "Ephb2", "FAPs", 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 3, 1, 3.821776, 11.465328, 0.8068087787009701, 0.8068087787009701, 14.89155655362667, 134.02400898264, 0.3241116275162939, 0.3241116275162939),
    @("ECs", "Efnb1", "Ephb2", "sCs", 3, 1, 3.896501666666667, 11.689505, 0.401720501899026, 0.401720501899026, 3, 1, 0.9091623333333333, 2.727487, 0.1919317489558758, 0.1919317489558758, 3.542552547103889, 31.882972923935, 0.07710291852091226, 0.07710291852091226),
    @("FAPs", "Efnb1", "Ephb2", "ECs", 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 1, 0.3333333333333333, 0.005966, 0.017898, 0.00125947234315407, 0.00125947234315407, 0.01931337384, 0.17382036456, 0.0004203515599978392, 0.0004203515599978392),
    @("FAPs", "Efnb1", "Ephb2", "FAPs", 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 3, 1, 3.821776, 11.465328, 0.8068087787009701, 0.8068087787009701, 12.37200613824, 111.34805524416, 0.2692741373721592, 0.2692741373721592),
    @("FAPs", "Efnb1", "Ephb2", "sCs", 3, 1, 3.23724, 9.71172, 0.3337521163387849, 0.3337521163387849, 3, 1, 0.9091623333333333, 2.727487, 0.1919317489558758, 0.1919317489558758, 2.94317667196, 26.48859004764, 0.0640576274066279, 0.0640576274066279),
    @("sCs", "Efnb1", "Ephb2", "ECs", 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 1, 0.3333333333333333, 0.005966, 0.017898, 0.00125947234315407, 0.00125947234315407, 0.01530751706066667, 0.137767653546, 0.0003331649213364355, 0.0003331649213364355),
    @("sCs", "Efnb1", "Ephb2", "FAPs", 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 3, 1, 3.821776, 11.465328, 0.8068087787009701, 0.8068087787009701, 9.805883560517334, 88.252952044656, 0.2134230138125171, 0.2134230138125171),
    @("sCs", "Efnb1", "Ephb2", "sCs", 3, 1, 2.565792333333333, 7.697377, 0.2645273817621892, 0.2645273817621892, 3, 1, 0.9091623333333333, 2.727487, 0.1919317489558758, 0.1919317489558758, 2.332721744622111, 20.994495701599, 0.05077120302833559, 0.05077120302833559)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $values = $rowData[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $values[$j]
    }
}
